# Weekly update: a new price record (week) was added to the data set.
# This inserts one new row at row 132 (pushing the existing rows 132..235
# down to 133..236) and populates the new row with the latest record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 132; this shifts rows 132..235
# down to 133..236, carrying all of their existing data/formatting with
# them (matches the diff, where every row from 133 downward now holds the
# content that used to belong to the row above it, and the final row 236
# holds what used to be row 235's content).
$ws.Rows(132).Insert()

# Populate the freshly inserted row 132 with the new weekly record.
$ws.Range("A132").Value = 11
$ws.Range("B132").Value = "Vega Monumental Concepción"
$ws.Range("C132").Value = "Bíobío"
$ws.Range("D132").Value = 45090
$ws.Range("E132").Value = 8
$ws.Range("F132").Value = 100112032
$ws.Range("G132").Value = "Zapallo italiano"
$ws.Range("H132").Value = "Sin especificar"
$ws.Range("I132").Value = "Primera"
$ws.Range("J132").Value = 100
$ws.Range("K132").Value = 11000
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = 11500
$ws.Range("N132").Value = "$/caja 50 unidades"
$ws.Range("O132").Value = "Región de Arica y Parinacota"
$ws.Range("P132").Value = 230
$ws.Range("Q132").Value = 50
$ws.Range("R132").Value = "Hortaliza"
